function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "311.82"
Set-TextValue $ws.Range("E2") "1.13%"
Set-TextValue $ws.Range("D3") "39.30"
Set-TextValue $ws.Range("E3") "1.63%"
Set-TextValue $ws.Range("D4") "5.137"
Set-TextValue $ws.Range("E4") "0.85%"
Set-TextValue $ws.Range("D5") "0.08136"
Set-TextValue $ws.Range("E5") "0.05%"
Set-TextValue $ws.Range("D6") "1.981"
Set-TextValue $ws.Range("E6") "1.25%"
Set-TextValue $ws.Range("D7") "4.234"
Set-TextValue $ws.Range("E7") "0.75%"
Set-TextValue $ws.Range("D8") "8.132"
Set-TextValue $ws.Range("E8") "2.48%"
Set-TextValue $ws.Range("D9") "0.9276"
Set-TextValue $ws.Range("E9") "-0.14%"
Set-TextValue $ws.Range("D10") "0.1400"
Set-TextValue $ws.Range("E10") "-3.11%"
Set-TextValue $ws.Range("D11") "0.1930"
Set-TextValue $ws.Range("E11") "-1.41%"
Set-TextValue $ws.Range("D12") "0.09054"
Set-TextValue $ws.Range("E12") "-0.58%"
Set-TextValue $ws.Range("D13") "0.03516"
Set-TextValue $ws.Range("E13") "0.29%"
Set-TextValue $ws.Range("D14") "0.09815"
Set-TextValue $ws.Range("E14") "-0.05%"
Set-TextValue $ws.Range("D15") "0.001402"
Set-TextValue $ws.Range("E15") "-0.25%"
Set-TextValue $ws.Range("D16") "0.006012"
Set-TextValue $ws.Range("E16") "1.03%"
Set-TextValue $ws.Range("D17") "3.683"
Set-TextValue $ws.Range("E17") "1.53%"
Set-TextValue $ws.Range("E18") "-2.41%"
Set-TextValue $ws.Range("D19") "0.3455"
Set-TextValue $ws.Range("E19") "0.25%"
Set-TextValue $ws.Range("E20") "-1.68%"
Set-TextValue $ws.Range("D21") "4.651"
Set-TextValue $ws.Range("E21") "-3.63%"
Set-TextValue $ws.Range("D22") "0.2424"
Set-TextValue $ws.Range("E22") "0.86%"
Set-TextValue $ws.Range("E23") "-1.66%"
Set-TextValue $ws.Range("E24") "0.15%"
Set-TextValue $ws.Range("D25") "0.004802"
Set-TextValue $ws.Range("E25") "-0.91%"
Set-TextValue $ws.Range("E26") "-0.12%"
Set-TextValue $ws.Range("D27") "0.0004000"
Set-TextValue $ws.Range("E27") "-10.06%"
Set-TextValue $ws.Range("D39") "0.02124"
Set-TextValue $ws.Range("E39") "1.31%"
Set-TextValue $ws.Range("D40") "0.05193"
Set-TextValue $ws.Range("E40") "1.42%"
Set-TextValue $ws.Range("D41") "0.007437"
Set-TextValue $ws.Range("E41") "-0.56%"
Set-TextValue $ws.Range("D42") "0.009832"
Set-TextValue $ws.Range("E42") "-3.05%"
Set-TextValue $ws.Range("D43") "0.1366"
Set-TextValue $ws.Range("E43") "0.25%"
Set-TextValue $ws.Range("D44") "0.002130"
Set-TextValue $ws.Range("E44") "-0.59%"
Set-TextValue $ws.Range("D45") "0.009005"
Set-TextValue $ws.Range("E45") "-14.12%"
Set-TextValue $ws.Range("D46") "0.00006402"
Set-TextValue $ws.Range("E46") "2.71%"
Set-TextValue $ws.Range("E47") "-0.17%"
Set-TextValue $ws.Range("E48") "-37.60%"
Set-TextValue $ws.Range("D49") "0.002568"
Set-TextValue $ws.Range("E49") "-16.09%"
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "-0.17%"
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.17%"
